# cards-test.xlsx: keep only the "cardNumber" column, dropping the old
# "first" / "middle" / "last" header columns (B and C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns B and C entirely (the "middle" and "last" headers),
# shifting everything after them left - nothing remains to their right.
$ws.Range("B1:C1").EntireColumn.Delete()

# Rename the remaining header in A1 from "first" to "cardNumber".
$ws.Range("A1").Value = "cardNumber"

# Re-fit column A to the new, longer header text.
$ws.Columns("A:A").AutoFit()

# Leave the selection where it lands after the column deletion (B1).
$ws.Range("B1").Select()
